# Updated TestData for Portugal Market
#
# 1) Duplicate the "Swiss" sheet (closest template: same layout/styles) to
#    create a new "Portugal" sheet placed right after it.
# 2) Fill in the Portugal-specific ticket number and market name (written in
#    the same order the diff's sharedStrings.xml shows them appended so the
#    shared-string table indices line up: ticket number first, market name
#    second).
# 3) Resize the new sheet's columns / the three wrapped description rows to
#    match the narrower layout used for Portugal.
# 4) Fix up selections: Germany's selection collapses to the full used range,
#    and the newly added Portugal tab ends up the active/selected tab.

$wb = $excel.ActiveWorkbook

$swiss = $wb.Worksheets.Item("Swiss")
$swiss.Copy($null, $swiss)
$portugal = $wb.Worksheets.Item($swiss.Index + 1)
$portugal.Name = "Portugal"

# Ticket number (B4) then market name (B2) -- matches the order the two new
# shared strings were appended in the target workbook.
$portugal.Range("B4").Value = "NGC-3479/T3493"
$portugal.Range("B2").Value = "Portugal Market"

# Narrower column layout for the Portugal sheet.
$portugal.Columns("A").ColumnWidth = 24.833333333333336
$portugal.Columns("B").ColumnWidth = 15
$portugal.Columns("C").ColumnWidth = 11
$portugal.Columns("D").ColumnWidth = 11.666666666666666

# With the narrower column D, the wrapped description cells spill onto a
# second line, so those rows end up twice the default height.
$portugal.Rows("3").RowHeight = 28.8
$portugal.Rows("4").RowHeight = 28.8
$portugal.Rows("5").RowHeight = 28.8

$portugal.Range("B2").Select() | Out-Null

# Germany's selection collapses from "A8:A15" to the whole used range.
$germany = $wb.Worksheets.Item("Germany")
$germany.Range("A1:D15").Select() | Out-Null

# Leave Portugal as the active sheet/tab.
$portugal.Activate() | Out-Null
